$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "capture cost" column (F) values per the new point progression.
# Old groups (row ranges -> old value -> new value):
#   F5:F11   2 -> 6
#   F12:F18  3 -> 7
#   F19:F25  4 -> 2
#   F26:F32  5 -> 3
#   F33:F39  6 -> 5
#   F40:F46  7 -> 4
$ws.Range("F5:F11").Value = 6
$ws.Range("F12:F18").Value = 7
$ws.Range("F19:F25").Value = 2
$ws.Range("F26:F32").Value = 3
$ws.Range("F33:F39").Value = 5
$ws.Range("F40:F46").Value = 4

# Update the active selection to match the author's final cursor position.
$ws.Range("F50").Select()
